$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 0.04224063622150328
$ws.Range("G3").Value = 0.04623746702836513
$ws.Range("G4").Value = -0.4471146569619693
$ws.Range("G5").Value = -0.4764070097706937
$ws.Range("G6").Value = 0.2308374342167992
$ws.Range("G7").Value = 0.2575500040655732
$ws.Range("G8").Value = 0.180416794454186
$ws.Range("G9").Value = 0.1807443943935199
$ws.Range("G10").Value = -0.005010809582411317
$ws.Range("G11").Value = -0.01434303962103854
$ws.Range("G12").Value = 0.1270353510948031
$ws.Range("G13").Value = 0.1482012788414959
$ws.Range("G14").Value = 0.259775015066963
$ws.Range("G15").Value = 0.2757615043332422
$ws.Range("G16").Value = 0.147924612147385
$ws.Range("G17").Value = 0.1504447085489411
$ws.Range("G18").Value = -0.007563514362421058
$ws.Range("G19").Value = -0.009686113988481531
$ws.Range("G20").Value = 0.1437065336997196
$ws.Range("G21").Value = 0.1421860672504419
$ws.Range("G22").Value = 0.1755136728354508
$ws.Range("G23").Value = 0.1880654429515226
$ws.Range("G24").Value = -0.1067871885409982
$ws.Range("G25").Value = -0.09419123910351263
$ws.Range("G26").Value = 0.2253044569855425
$ws.Range("G27").Value = 0.2473133869698949
$ws.Range("G28").Value = 0.0664943628212393
$ws.Range("G29").Value = 0.09116362591383118
$ws.Range("H2").Value = -12.41946449679861
$ws.Range("H3").Value = 20.54630192738958
$ws.Range("H4").Value = 1.248615549802442
$ws.Range("H5").Value = 0.5070575404202432
$ws.Range("H6").Value = -1.195992318922308
$ws.Range("H7").Value = 16.76321321471805
$ws.Range("H8").Value = 8.159371082553912
$ws.Range("H9").Value = 5.077028734786327
$ws.Range("H10").Value = -5.724839673976486
$ws.Range("H11").Value = 2.130670431367979
$ws.Range("H12").Value = -7.085110225164216
$ws.Range("H13").Value = 18.90577166819024
$ws.Range("H14").Value = 5.030269125358084
$ws.Range("H15").Value = 9.135969428343724
$ws.Range("H16").Value = -3.612878830131451
$ws.Range("H17").Value = -0.3777698276636051
$ws.Range("H18").Value = 53.81033488276903
$ws.Range("H19").Value = -1050.282518685604
$ws.Range("H20").Value = 3.637878091552599
$ws.Range("H21").Value = -0.634830897505152
$ws.Range("H22").Value = -5.747966597357117
$ws.Range("H23").Value = 4.794892760208906
$ws.Range("H24").Value = -13.1253020509411
$ws.Range("H25").Value = 5.438509701481843
$ws.Range("H26").Value = -2.101319342488557
$ws.Range("H27").Value = 6.338309698990314
$ws.Range("H28").Value = 13.07937654360114
$ws.Range("H29").Value = 29.15526533623171
$ws.Range("I2").Value = -32.09104132959953
